# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.956.01'
$ws.Range('E2').Value = '  -0.75%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.618.89'
$ws.Range('E3').Value = '  -1.16%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.10'
$ws.Range('E5').Value = '  -1.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.511'
$ws.Range('E6').Value = '  -1.30%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.250'
$ws.Range('E8').Value = '  -1.48%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0625'
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.90'
$ws.Range('E10').Value = '  -0.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0838'
$ws.Range('E11').Value = '  -1.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.846.84'
$ws.Range('E12').Value = '  -1.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.616.64'
$ws.Range('E13').Value = '  -1.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.09'
$ws.Range('E14').Value = '  -0.97%  '
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.955.10'
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.11'
$ws.Range('E17').Value = '  -3.43%  '
$ws.Range('E18').Value = '  -0.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '212.73'
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.78'
$ws.Range('E21').Value = '  -1.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.32'
$ws.Range('E22').Value = '  -2.21%  '
$ws.Range('E23').Value = '  -8.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.95'
$ws.Range('E24').Value = '  -2.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.63'
$ws.Range('E25').Value = '  -0.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.46'
$ws.Range('E26').Value = '  +1.07%  '
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('E28').Value = '  -3.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.44'
$ws.Range('E29').Value = '  -1.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0506'
$ws.Range('E30').Value = '  -0.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.16'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.28'
$ws.Range('E32').Value = '  -2.80%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.700'
$ws.Range('E33').Value = '  +27.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.98'
$ws.Range('E34').Value = '  -1.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.334.66'
$ws.Range('E35').Value = '  +2.69%  '
$ws.Range('E36').Value = '  -1.45%  '
$ws.Range('E37').Value = '  -0.68%  '
$ws.Range('E38').Value = '  -0.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.836'
$ws.Range('E39').Value = '  -2.09%  '
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.21'
$ws.Range('E41').Value = '  -2.12%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.794'
$ws.Range('E42').Value = '  -2.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.32'
$ws.Range('E43').Value = '  -0.43%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.63'
$ws.Range('E44').Value = '  +1.87%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.757.88'
$ws.Range('E45').Value = '  -1.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.68'
$ws.Range('E46').Value = '  -1.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.62'
$ws.Range('E47').Value = '  +1.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.794'
$ws.Range('E48').Value = '  +6.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0515'
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('E50').Value = '  +2.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.54'
$ws.Range('E51').Value = '  -0.89%  '
